$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Insert a new worksheet for "2022-Q4" right before the current 2nd sheet
# (today that's "2021-Q1"), so the final tab order becomes:
# 总计, 2022-Q4, 2021-Q1, 2020-Q4
$newSheet = $wb.Worksheets.Add($wb.Worksheets.Item(2))
$newSheet.Name = "2022-Q4"

# ---- Header row (B1:H1) ------------------------------------------------
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Match the bold/bordered/centered header style used elsewhere in the workbook
$ws1.Range("B1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122) # xlPasteFormats

# ---- Row 2: fund 004618 -------------------------------------------------
$newSheet.Range("A2").Value = 0
$ws1.Range("A2").Copy()
$newSheet.Range("A2").PasteSpecial(-4122)

$newSheet.Range("B2:G2").NumberFormat = "@"
$newSheet.Range("B2").Value = "004618"
$newSheet.Range("C2").Value = "建信鑫稳回报灵活配置混合C"
$newSheet.Range("D2").Value = "0.75"
$newSheet.Range("E2").Value = "26.14"
$newSheet.Range("F2").Value = "0.63"
$newSheet.Range("G2").Value = "0.0047"
$newSheet.Range("H2").Value = 7

# Restore the default (General) format now that the text is stored, so the
# cells end up styled the same as a plain untouched cell
$newSheet.Range("A1").Copy()
$newSheet.Range("B2:G2").PasteSpecial(-4122)

# ---- Row 3: fund 004617 -------------------------------------------------
$newSheet.Range("A3").Value = 1
$ws1.Range("A2").Copy()
$newSheet.Range("A3").PasteSpecial(-4122)

$newSheet.Range("B3:G3").NumberFormat = "@"
$newSheet.Range("B3").Value = "004617"
$newSheet.Range("C3").Value = "建信鑫稳回报灵活配置混合A"
$newSheet.Range("D3").Value = "0.31"
$newSheet.Range("E3").Value = "26.14"
$newSheet.Range("F3").Value = "0.63"
$newSheet.Range("G3").Value = "0.0020"
$newSheet.Range("H3").Value = 7

$newSheet.Range("A1").Copy()
$newSheet.Range("B3:G3").PasteSpecial(-4122)

# ---- Update the "总计" (totals) sheet -----------------------------------
# Shift existing rows down to make room for the new 2022-Q4 entry on top,
# then append 2020-Q4 as the new row 4.
$ws1.Range("A3").Copy()
$ws1.Range("A4").PasteSpecial(-4122) # copy A3's style onto the new A4

$ws1.Range("B2").Value = "2022-Q4"
$ws1.Range("C2").Value = 2
$ws1.Range("D2").Value = 0.01

$ws1.Range("B3").Value = "2021-Q1"
$ws1.Range("C3").Value = 1
$ws1.Range("D3").Value = 0.61

$ws1.Range("A4").Value = 2
$ws1.Range("B4").Value = "2020-Q4"
$ws1.Range("C4").Value = 1
$ws1.Range("D4").Value = 0.54

# Keep "2020-Q4" as the active/selected tab, same as before the edit
# (adding a worksheet makes it active by default).
$wb.Worksheets.Item("2020-Q4").Activate()
